$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = 29
$ws.Cells.Item(2, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(2, 3).Value = 0.9340000152587891
$ws.Cells.Item(2, 4).Value = 0.0009999871253967285
$ws.Cells.Item(2, 5).Value = 1

$ws.Cells.Item(3, 1).Value = 21
$ws.Cells.Item(3, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(3, 3).Value = 0.9329999685287476
$ws.Cells.Item(3, 4).Value = 0.0009999871253967285
$ws.Cells.Item(3, 5).Value = 2

$ws.Cells.Item(4, 1).Value = 31
$ws.Cells.Item(4, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(4, 3).Value = 0.9300000369548798
$ws.Cells.Item(4, 4).Value = 0.00600007176399231
$ws.Cells.Item(4, 5).Value = 3

$ws.Cells.Item(5, 1).Value = 18
$ws.Cells.Item(5, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(5, 3).Value = 0.9294999837875366
$ws.Cells.Item(5, 4).Value = 0.0004999637603759766
$ws.Cells.Item(5, 5).Value = 4

$ws.Cells.Item(6, 1).Value = 23
$ws.Cells.Item(6, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(6, 3).Value = 0.9284999966621399
$ws.Cells.Item(6, 4).Value = 0.003499984741210938
$ws.Cells.Item(6, 5).Value = 5

$ws.Cells.Item(7, 1).Value = 22
$ws.Cells.Item(7, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(7, 3).Value = 0.9280000627040863
$ws.Cells.Item(7, 4).Value = 0.001000016927719116
$ws.Cells.Item(7, 5).Value = 6

$ws.Cells.Item(8, 1).Value = 14
$ws.Cells.Item(8, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(8, 3).Value = 0.9275000393390656
$ws.Cells.Item(8, 4).Value = 0.005499929189682007
$ws.Cells.Item(8, 5).Value = 7

$ws.Cells.Item(9, 1).Value = 16
$ws.Cells.Item(9, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(9, 3).Value = 0.9255000054836273
$ws.Cells.Item(9, 4).Value = 0.008500009775161743
$ws.Cells.Item(9, 5).Value = 8

$ws.Cells.Item(10, 1).Value = 4
$ws.Cells.Item(10, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(10, 3).Value = 0.9254999756813049
$ws.Cells.Item(10, 4).Value = 0.004499971866607666
$ws.Cells.Item(10, 5).Value = 9

$ws.Cells.Item(11, 1).Value = 28
$ws.Cells.Item(11, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(11, 3).Value = 0.9245000183582306
$ws.Cells.Item(11, 4).Value = 0.003500014543533325
$ws.Cells.Item(11, 5).Value = 10

$ws.Cells.Item(12, 1).Value = 30
$ws.Cells.Item(12, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(12, 3).Value = 0.9245000183582306
$ws.Cells.Item(12, 4).Value = 0.006499975919723511
$ws.Cells.Item(12, 5).Value = 10

$ws.Cells.Item(13, 1).Value = 5
$ws.Cells.Item(13, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(13, 3).Value = 0.9239999353885651
$ws.Cells.Item(13, 4).Value = 0.00600007176399231
$ws.Cells.Item(13, 5).Value = 12

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(14, 3).Value = 0.92249995470047
$ws.Cells.Item(14, 4).Value = 0.006500065326690674
$ws.Cells.Item(14, 5).Value = 13

$ws.Cells.Item(15, 1).Value = 20
$ws.Cells.Item(15, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(15, 3).Value = 0.918999969959259
$ws.Cells.Item(15, 4).Value = 0.01100003719329834
$ws.Cells.Item(15, 5).Value = 14

$ws.Cells.Item(16, 1).Value = 6
$ws.Cells.Item(16, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(16, 3).Value = 0.9174999892711639
$ws.Cells.Item(16, 4).Value = 0.01250001788139343
$ws.Cells.Item(16, 5).Value = 15

$ws.Cells.Item(17, 1).Value = 7
$ws.Cells.Item(17, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(17, 3).Value = 0.9149999916553497
$ws.Cells.Item(17, 4).Value = 0.009000033140182495
$ws.Cells.Item(17, 5).Value = 16

$ws.Cells.Item(18, 1).Value = 2
$ws.Cells.Item(18, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(18, 3).Value = 0.9095000624656677
$ws.Cells.Item(18, 4).Value = 0.01049995422363281
$ws.Cells.Item(18, 5).Value = 17

$ws.Cells.Item(19, 1).Value = 25
$ws.Cells.Item(19, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(19, 3).Value = 0.9085000157356262
$ws.Cells.Item(19, 4).Value = 0.002500057220458984
$ws.Cells.Item(19, 5).Value = 18

$ws.Cells.Item(20, 1).Value = 24
$ws.Cells.Item(20, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(20, 3).Value = 0.9075000286102295
$ws.Cells.Item(20, 4).Value = 0.02150005102157593
$ws.Cells.Item(20, 5).Value = 19

$ws.Cells.Item(21, 1).Value = 26
$ws.Cells.Item(21, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(21, 3).Value = 0.9074999690055847
$ws.Cells.Item(21, 4).Value = 0.01750004291534424
$ws.Cells.Item(21, 5).Value = 20

$ws.Cells.Item(22, 1).Value = 17
$ws.Cells.Item(22, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(22, 3).Value = 0.9049999117851257
$ws.Cells.Item(22, 4).Value = 0.005999982357025146
$ws.Cells.Item(22, 5).Value = 21

$ws.Cells.Item(23, 1).Value = 27
$ws.Cells.Item(23, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(23, 3).Value = 0.9015000462532043
$ws.Cells.Item(23, 4).Value = 0.00850003957748413
$ws.Cells.Item(23, 5).Value = 22

$ws.Cells.Item(24, 1).Value = 19
$ws.Cells.Item(24, 2).Value = '{''anOptimizer'': ''adam'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(24, 3).Value = 0.8959999978542328
$ws.Cells.Item(24, 4).Value = 0.003999978303909302
$ws.Cells.Item(24, 5).Value = 23

$ws.Cells.Item(25, 1).Value = 12
$ws.Cells.Item(25, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(25, 3).Value = 0.8950000107288361
$ws.Cells.Item(25, 4).Value = 0.03399994969367981
$ws.Cells.Item(25, 5).Value = 24

$ws.Cells.Item(26, 1).Value = 15
$ws.Cells.Item(26, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 10, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(26, 3).Value = 0.8944999575614929
$ws.Cells.Item(26, 4).Value = 0.007499992847442627
$ws.Cells.Item(26, 5).Value = 25

$ws.Cells.Item(27, 1).Value = 10
$ws.Cells.Item(27, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''softmax''}'
$ws.Cells.Item(27, 3).Value = 0.8794999718666077
$ws.Cells.Item(27, 4).Value = 0.003499984741210938
$ws.Cells.Item(27, 5).Value = 26

$ws.Cells.Item(28, 1).Value = 3
$ws.Cells.Item(28, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(28, 3).Value = 0.8790000081062317
$ws.Cells.Item(28, 4).Value = 0.01899999380111694
$ws.Cells.Item(28, 5).Value = 27

$ws.Cells.Item(29, 1).Value = 11
$ws.Cells.Item(29, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 128, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(29, 3).Value = 0.8719999492168427
$ws.Cells.Item(29, 4).Value = 0.0379999577999115
$ws.Cells.Item(29, 5).Value = 28

$ws.Cells.Item(30, 1).Value = 1
$ws.Cells.Item(30, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(30, 3).Value = 0.8509999513626099
$ws.Cells.Item(30, 4).Value = 0.00700002908706665
$ws.Cells.Item(30, 5).Value = 29

$ws.Cells.Item(31, 1).Value = 8
$ws.Cells.Item(31, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(31, 3).Value = 0.8479999899864197
$ws.Cells.Item(31, 4).Value = 0.05999994277954102
$ws.Cells.Item(31, 5).Value = 30

$ws.Cells.Item(32, 1).Value = 9
$ws.Cells.Item(32, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 20, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''sigmoid''}'
$ws.Cells.Item(32, 3).Value = 0.8449999988079071
$ws.Cells.Item(32, 4).Value = 0.01299998164176941
$ws.Cells.Item(32, 5).Value = 31

$ws.Cells.Item(33, 1).Value = 0
$ws.Cells.Item(33, 2).Value = '{''anOptimizer'': ''rmsprop'', ''batch_size'': 25, ''epochs'': 5, ''hidUnit'': 256, ''outActivation'': ''softmax''}'
$ws.Cells.Item(33, 3).Value = 0.8259999752044678
$ws.Cells.Item(33, 4).Value = 0.04900002479553223
$ws.Cells.Item(33, 5).Value = 32
